# Regenerate the "K" column (column G) values for the save_data sheet.
# These values replace the previous "Strike#" derived numbers with the
# recalculated K counts (std/mean regen + s_vals write-back).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 0
    8  = 0
    9  = 2
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 2
    16 = 0
    17 = 0
    18 = 2
    19 = 2
    20 = 0
    21 = 0
    22 = 2
    23 = 0
    24 = 1
    25 = 0
    27 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
